$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.643.20"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.454.07"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.68"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.94"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.530"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.111"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.348"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.58"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.885.91"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.518.90"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.457.72"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.63"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -6.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.75"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "320.86"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.13"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.21"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.91"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.06"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "643.87"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.559.95"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0956"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.995"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.41"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.48%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.81"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.25%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.42%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.64"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "150.41"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.365"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.50"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.32"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.71"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.71"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₆0309"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.38%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "152.59"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.38"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.54"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.604"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.12"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0503"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.60%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.25%  "

Write-Host "Applied cryptos update"